$wb = $excel.ActiveWorkbook

# --- Sheet 1: Psychosis_medications ---
$ws1 = $wb.Worksheets.Item("Psychosis_medications")

# Header renames
$ws1.Range("G1").Value = "Inventor_manufacturer_first"
$ws1.Range("H1").Value = "Indications_current"

# Row 3 (Reserpine): the drug class value ("Rauwolfia alkaloid") had been
# mistakenly placed in the Treatment_class column (C3) while Drug_class
# (D3) was empty. Fix: give Reserpine its Treatment_class and move the
# drug class into D3.
$ws1.Range("C3").Value = "First-generation antipsychotic"
$ws1.Range("D3").Value = "Rauwolfia alkaloid"

# --- Sheet 2: Metadata_psychosis_medications ---
$ws2 = $wb.Worksheets.Item("Metadata_psychosis_medications")

# Row 8 documents the Inventor_manufacturer(_first) column.
$ws2.Range("A8").Value = "Inventor_manufacturer_first"
$ws2.Range("B8").Value = "Name of the scientist, group, company, or institution which is credited with its invention; if unavailable, its first manufacturer"
